# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date serial in A1 by one day (2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: update the price list (column D, rows 32-39)
$ws.Range("D32").Value = 219.124
$ws.Range("D33").Value = 313.543
$ws.Range("D34").Value = 417.457
$ws.Range("D35").Value = 429.794
$ws.Range("D36").Value = 563.266
$ws.Range("D37").Value = 644.069
$ws.Range("D38").Value = 771.2670000000001
$ws.Range("D39").Value = 918.41
